$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 98.5
$ws.Range("J9").Value = 99
$ws.Range("L9").Value = 99
$ws.Range("N9").Value = -437
$ws.Range("H12").Value = 1329.4
$ws.Range("I12").Value = 1199
$ws.Range("J12").Value = 1416.3334
$ws.Range("K12").Value = 1199
$ws.Range("L12").Value = 1416.3334
$ws.Range("M12").Value = -1029
$ws.Range("N12").Value = -1756.3334
$ws.Range("H15").Value = 428.7
$ws.Range("I15").Value = 428.7
$ws.Range("K15").Value = 1286.1
$ws.Range("M15").Value = -1117.1
$ws.Range("H92").Value = 611.6429000000001
$ws.Range("I92").Value = 385.1111
$ws.Range("J92").Value = 1019.4
$ws.Range("K92").Value = 385.1111
$ws.Range("L92").Value = 1019.4
$ws.Range("M92").Value = 862.8888999999999
$ws.Range("N92").Value = -3515.4
$ws.Range("H108").Value = 44400
$ws.Range("J108").Value = 44400
$ws.Range("L108").Value = 44400
$ws.Range("N108").Value = -52080
$ws.Range("H138").Value = 3723.8657
$ws.Range("I138").Value = 4859.375
$ws.Range("J138").Value = 3367.6274
$ws.Range("K138").Value = 14578.125
$ws.Range("L138").Value = 10102.8822
$ws.Range("M138").Value = -9438.125
$ws.Range("N138").Value = -20382.8822
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14103.8
$ws.Range("I32").Value = 11097.17
$ws.Range("J32").Value = 25190.75
$ws.Range("K32").Value = 11097.17
$ws.Range("L32").Value = 25190.75
$ws.Range("M32").Value = -10810.17
$ws.Range("N32").Value = -25764.75
$ws.Range("H109").Value = 68744.5
$ws.Range("J109").Value = 68744.5
$ws.Range("L109").Value = 68744.5
$ws.Range("N109").Value = -71518.5
$ws.Range("H132").Value = 1681.7
$ws.Range("I132").Value = 1379
$ws.Range("J132").Value = 1967.5834
$ws.Range("K132").Value = 4137
$ws.Range("L132").Value = 5902.7502
$ws.Range("M132").Value = -1607
$ws.Range("N132").Value = -10962.7502
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 287428.16
$ws.Range("J86").Value = 2000000
$ws.Range("L86").Value = 2000000
$ws.Range("N86").Value = -2002246
$ws.Range("H89").Value = 287428.16
$ws.Range("J89").Value = 2000000
$ws.Range("L89").Value = 10000000
$ws.Range("N89").Value = -10011232
$ws.Range("H99").Value = 1675
$ws.Range("I99").Value = 1733.3334
$ws.Range("J99").Value = 1500
$ws.Range("K99").Value = 1733.3334
$ws.Range("L99").Value = 1500
$ws.Range("M99").Value = -235.3334
$ws.Range("N99").Value = -4496
$ws.Range("H134").Value = 3986.9812
$ws.Range("I134").Value = 3988.5557
$ws.Range("K134").Value = 11965.6671
$ws.Range("M134").Value = -9430.667099999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 65.25
$ws.Range("I7").Value = 65.25
$ws.Range("K7").Value = 65.25
$ws.Range("M7").Value = 47.75
$ws.Range("H99").Value = 8003
$ws.Range("I99").Value = 8670.666999999999
$ws.Range("K99").Value = 8670.666999999999
$ws.Range("M99").Value = -7172.666999999999
$ws.Range("H126").Value = 8003
$ws.Range("I126").Value = 8670.666999999999
$ws.Range("K126").Value = 26012.001
$ws.Range("M126").Value = -23542.001
$ws.Range("H132").Value = 2532.9375
$ws.Range("I132").Value = 1611.091
$ws.Range("K132").Value = 4833.272999999999
$ws.Range("M132").Value = -2303.272999999999
$ws.Range("H134").Value = 977.98114
$ws.Range("I134").Value = 815.61365
$ws.Range("J134").Value = 1771.7778
$ws.Range("K134").Value = 2446.84095
$ws.Range("L134").Value = 5315.3334
$ws.Range("M134").Value = 88.15905000000021
$ws.Range("N134").Value = -10385.3334
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 695.1429000000001
$ws.Range("I5").Value = 532.1
$ws.Range("J5").Value = 1102.75
$ws.Range("K5").Value = 1596.3
$ws.Range("L5").Value = 3308.25
$ws.Range("M5").Value = -1484.3
$ws.Range("N5").Value = -3532.25
$ws.Range("H7").Value = 1136.875
$ws.Range("I7").Value = 200
$ws.Range("J7").Value = 1270.7142
$ws.Range("K7").Value = 600
$ws.Range("L7").Value = 3812.1426
$ws.Range("M7").Value = -488
$ws.Range("N7").Value = -4036.1426
$ws.Range("H133").Value = 1500
$ws.Range("I133").Value = 1500
$ws.Range("K133").Value = 4500
$ws.Range("M133").Value = 560
$ws.Range("H135").Value = 695.1429000000001
$ws.Range("I135").Value = 532.1
$ws.Range("J135").Value = 1102.75
$ws.Range("K135").Value = 4788.900000000001
$ws.Range("L135").Value = 9924.75
$ws.Range("M135").Value = -2253.900000000001
$ws.Range("N135").Value = -14994.75
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2065.2
$ws.Range("I80").Value = 1831.5
$ws.Range("J80").Value = 3000
$ws.Range("K80").Value = 1831.5
$ws.Range("L80").Value = 3000
$ws.Range("M80").Value = -833.5
$ws.Range("N80").Value = -4996
$ws.Range("H83").Value = 2065.2
$ws.Range("I83").Value = 1831.5
$ws.Range("J83").Value = 3000
$ws.Range("K83").Value = 9157.5
$ws.Range("L83").Value = 15000
$ws.Range("M83").Value = -4165.5
$ws.Range("N83").Value = -24984
$ws.Range("H102").Value = 2788.5833
$ws.Range("I102").Value = 2632.875
$ws.Range("K102").Value = 2632.875
$ws.Range("M102").Value = -1010.875
$ws.Range("H113").Value = 1811
$ws.Range("I113").Value = 1811
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1811
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 359
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 634448.9399999999
$ws.Range("I132").Value = 898832.5
$ws.Range("J132").Value = 2865.9443
$ws.Range("K132").Value = 2696497.5
$ws.Range("L132").Value = 8597.832900000001
$ws.Range("M132").Value = -2693967.5
$ws.Range("N132").Value = -13657.8329
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3327.625
$ws.Range("I7").Value = 3283.8333
$ws.Range("K7").Value = 3283.8333
$ws.Range("M7").Value = -3171.8333
$ws.Range("H40").Value = 18958.584
$ws.Range("I40").Value = 29625.75
$ws.Range("J40").Value = 13625
$ws.Range("K40").Value = 29625.75
$ws.Range("L40").Value = 13625
$ws.Range("M40").Value = -29489.75
$ws.Range("N40").Value = -13897
$ws.Range("H122").Value = 4270.8887
$ws.Range("I122").Value = 3875.2222
$ws.Range("K122").Value = 11625.6666
$ws.Range("M122").Value = -9175.6666
$ws.Range("H126").Value = 3327.625
$ws.Range("I126").Value = 3283.8333
$ws.Range("K126").Value = 9851.499899999999
$ws.Range("M126").Value = -7381.499899999999
$ws.Range("H132").Value = 5354.205
$ws.Range("I132").Value = 3891.1853
$ws.Range("K132").Value = 11673.5559
$ws.Range("M132").Value = -9143.555899999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 11280
$ws.Range("J54").Value = 11280
$ws.Range("L54").Value = 11280
$ws.Range("N54").Value = -12320
$ws.Range("H122").Value = 30590.965
$ws.Range("J122").Value = 2825.5454
$ws.Range("L122").Value = 8476.636200000001
$ws.Range("N122").Value = -13376.6362
